$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games) for rows 2-5, columns B-G.
$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.002658071450198252
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.7273727411070765

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4.327115817150455

$ws.Range("B4").Value = 0.1169995834814548
$ws.Range("C4").Value = 0.3048912486333797
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1.104883657715537

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.1496068669990043
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.582307763322248
